# Append a literal "." right after "UNICEN-TUDAI" in its own run,
# matching the target diff which adds a brand new <w:r> (with the same
# run formatting: Times New Roman, sz/szCs 40, lang es-AR) containing
# only the new period character, instead of folding it into the
# existing "UNICEN-TUDAI" run.

$d = $word.ActiveDocument

# Locate the "UNICEN-TUDAI" text.
$rng = $d.Content
$found = $rng.Find.Execute("UNICEN-TUDAI", $false, $false, $false, $false, `
                            $false, $true, 1, $false, "", 0)

if ($found) {
    # Collapse the found range to its end (right after "UNICEN-TUDAI")
    # and insert the new period there.
    $rng.Collapse(0)
    $rng.InsertAfter(".")

    # The freshly inserted "." currently lives in the same run as
    # "UNICEN-TUDAI" (identical formatting gets merged). Touch a
    # formatting property on just the new character (set then restore)
    # so it is materialised as its own separate <w:r> run, matching
    # the target document structure, while keeping the final
    # formatting identical to the surrounding text.
    $newCharRange = $d.Range($rng.Start, $rng.Start + 1)
    $newCharRange.Font.Bold = 1
    $newCharRange.Font.Bold = 0
}
